$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bad Drivers table updates
$ws.Range("C3").Value = 655
$ws.Range("D3").Value = 84.5
$ws.Range("D4").Value = 98.90000000000001
$ws.Range("C5").Value = 717

# Good Drivers table updates
$ws.Range("B13").Value = 449371
$ws.Range("B14").Value = 77999
